# Update "want to go" counts (column F) on several sheets to match the
# latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1084
$ws1.Range("F4").Value = 1664
$ws1.Range("F6").Value = 170

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 12

# 全部类型 (All types) - aggregated view of the sheets above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1084
$ws4.Range("F4").Value = 1664
$ws4.Range("F5").Value = 12
$ws4.Range("F7").Value = 170
